$wb = $excel.ActiveWorkbook

# --- Update sheet "sets": D4 (home_points for set_id=3) changes from 7 to 8 ---
$wsSets = $wb.Worksheets.Item("sets")
$wsSets.Range("D4").Value = 8

# --- Append a new row (71) to sheet "rallies" ---
$wsRallies = $wb.Worksheets.Item("rallies")

$wsRallies.Cells.Item(71, 1).Value = 70      # A71 rally_id
$wsRallies.Cells.Item(71, 2).Value = 1       # B71 match_id
$wsRallies.Cells.Item(71, 3).Value = 3       # C71 set_number
$wsRallies.Cells.Item(71, 4).Value = 8       # D71 rally_no
$wsRallies.Cells.Item(71, 5).Value = "NOS"   # E71 side
$wsRallies.Cells.Item(71, 6).Value = "'"     # F71 position (empty text cell)
$wsRallies.Cells.Item(71, 6).Style = "Normal"
$wsRallies.Cells.Item(71, 7).Value = 2       # G71 player_number
$wsRallies.Cells.Item(71, 8).Value = "LINHA" # H71 action
$wsRallies.Cells.Item(71, 9).Value = "PONTO" # I71 result
$wsRallies.Cells.Item(71, 10).Value = "NOS"  # J71 who_scored
$wsRallies.Cells.Item(71, 11).Value = 8      # K71 score_home
$wsRallies.Cells.Item(71, 12).Value = 0      # L71 score_away
$wsRallies.Cells.Item(71, 13).Value = "1 2 l"    # M71 raw_text
$wsRallies.Cells.Item(71, 14).Value = "FRENTE"   # N71 position_zone
$wsRallies.Cells.Item(71, 15).Value = "FRENTE"   # O71 pos_fb
$wsRallies.Cells.Item(71, 16).Value = "FRENTE"   # P71 frente_fundo
